$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells AZ1:BE1 ---
$ws.Range("AZ1").Value = "Publicações JCR (total)"
$ws.Range("BA1").Value = "Publicações JCR > 1,5 (total)"
$ws.Range("BB1").Value = "Publicações JCR"
$ws.Range("BC1").Value = "Publicações JCR > 1,5"
$ws.Range("BD1").Value = "Aceitações JCR > 1,5"
$ws.Range("BE1").Value = "Artigos JCR > 1,5"

# Match the bold/centered/bordered header style used by the rest of row 1
$ws.Range("AY1").Copy()
$ws.Range("AZ1:BE1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- New data columns AZ2:BE40 ---
$newData = @(
    @(2,18,11,8,4,0,4),
    @(3,4,3,3,2,0,2),
    @(4,5,3,1,0,0,0),
    @(5,4,4,0,0,0,0),
    @(6,34,24,11,7,1,8),
    @(7,2,1,1,0,0,0),
    @(8,11,4,1,1,1,2),
    @(9,21,17,5,4,0,4),
    @(10,93,56,11,10,1,11),
    @(11,5,2,0,0,0,0),
    @(12,0,0,0,0,1,1),
    @(13,17,11,8,6,1,7),
    @(14,4,2,0,0,0,0),
    @(15,14,9,3,2,1,3),
    @(16,11,11,3,3,0,3),
    @(17,12,7,6,3,0,3),
    @(18,8,3,1,1,0,1),
    @(19,56,11,15,4,0,4),
    @(20,15,8,3,1,0,1),
    @(21,14,9,2,0,0,0),
    @(22,4,4,1,1,0,1),
    @(23,22,21,3,3,1,4),
    @(24,9,8,1,1,0,1),
    @(25,18,11,8,7,1,8),
    @(26,10,1,1,0,0,0),
    @(27,21,11,6,3,0,3),
    @(28,1,0,1,0,0,0),
    @(29,7,0,1,0,0,0),
    @(30,9,5,5,4,0,4),
    @(31,2,1,0,0,0,0),
    @(32,37,30,14,12,2,14),
    @(33,11,9,5,4,1,5),
    @(34,24,21,3,2,0,2),
    @(35,7,2,0,0,0,0),
    @(36,12,11,3,3,0,3),
    @(37,14,2,10,1,0,1),
    @(38,12,6,4,3,0,3),
    @(39,8,7,1,1,0,1),
    @(40,14,6,5,4,0,4)
)

foreach ($entry in $newData) {
    $r = $entry[0]
    $ws.Range("AZ" + $r).Value = $entry[1]
    $ws.Range("BA" + $r).Value = $entry[2]
    $ws.Range("BB" + $r).Value = $entry[3]
    $ws.Range("BC" + $r).Value = $entry[4]
    $ws.Range("BD" + $r).Value = $entry[5]
    $ws.Range("BE" + $r).Value = $entry[6]
}

# --- Updated existing values (AF/AG/AH/AV/AW) ---
$cellUpdates = @(
    @(2,"AF",762),
    @(2,"AH",283),
    @(2,"AV",42.33333333333334),
    @(6,"AF",2425),
    @(6,"AH",738),
    @(6,"AV",80.83333333333333),
    @(9,"AF",2113),
    @(9,"AH",492),
    @(9,"AV",117.3888888888889),
    @(10,"AF",9121),
    @(10,"AH",1526),
    @(10,"AV",260.6),
    @(11,"AF",676),
    @(11,"AG",15),
    @(11,"AH",95),
    @(11,"AV",39.76470588235294),
    @(11,"AW",0.8823529411764706),
    @(13,"AF",1975),
    @(13,"AH",805),
    @(13,"AV",329.1666666666667),
    @(15,"AF",1557),
    @(15,"AH",455),
    @(15,"AV",103.8),
    @(17,"AF",1502),
    @(17,"AH",497),
    @(17,"AV",107.2857142857143),
    @(18,"AF",634),
    @(18,"AH",87),
    @(18,"AV",28.81818181818182),
    @(19,"AF",1165),
    @(19,"AH",356),
    @(19,"AV",58.25),
    @(20,"AF",1230),
    @(20,"AH",303),
    @(20,"AV",136.6666666666667),
    @(22,"AF",449),
    @(22,"AH",190),
    @(22,"AV",49.88888888888889),
    @(23,"AF",2200),
    @(23,"AH",773),
    @(23,"AV",100),
    @(24,"AF",752),
    @(24,"AH",237),
    @(24,"AV",94),
    @(25,"AF",1803),
    @(25,"AH",483),
    @(25,"AV",150.25),
    @(27,"AF",1145),
    @(27,"AV",47.70833333333334),
    @(32,"AH",993),
    @(33,"AF",224),
    @(33,"AH",111),
    @(33,"AV",24.88888888888889),
    @(34,"AF",1422),
    @(34,"AH",343),
    @(34,"AV",40.62857142857143),
    @(35,"AF",1146),
    @(35,"AH",179),
    @(35,"AV",33.70588235294117),
    @(36,"AF",1126),
    @(36,"AH",224),
    @(36,"AV",59.26315789473684),
    @(37,"AF",135),
    @(37,"AH",102),
    @(37,"AV",33.75),
    @(38,"AF",1281),
    @(38,"AH",318),
    @(38,"AV",91.5),
    @(39,"AF",1208),
    @(39,"AV",86.28571428571429)
)

foreach ($u in $cellUpdates) {
    $ws.Range($u[1] + $u[0]).Value = $u[2]
}
